$d = $word.ActiveDocument

# Locate the paragraph containing "LOM3057" (the requirement line that
# precedes the blank paragraph and the "Ver no Jupiter..." / "© 2020..."
# footer paragraphs that must be removed).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "LOM3057") {
        $target = $i
        break
    }
}

# The three paragraphs to delete are the ones immediately following the
# "LOM3057..." paragraph: an empty "Normal" paragraph, the
# "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, and the
# "© 2020 . Contact: ..." paragraph.
$firstToDelete = $target + 1
$lastToDelete = $target + 3

$startPara = $d.Paragraphs.Item($firstToDelete)
$endPara = $d.Paragraphs.Item($lastToDelete)

$r = $d.Range($startPara.Range.Start, $endPara.Range.End)
$r.Delete()
